$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new "Multi" worksheet at the very end of the workbook (after
#    "Variance").
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMulti = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsMulti.Name = "Multi"
$wsMulti.Activate()

# ---------------------------------------------------------------------------
# 2. Build the new "Multi" sheet content FIRST: a wider, multi-currency /
#    multi-milestone cash-flow table. Doing this before touching the other
#    sheets keeps the shared-string table insertion order (Unit_ID,
#    Milestone_ID, Currency, FX_Reference, Unit_A, MS1, USD, USD_THB)
#    matching the source order of authorship.
# ---------------------------------------------------------------------------
$headers = @("Date","Unit_ID","Milestone_ID","Currency","FX_Reference","Planned Inflow","Planned Outflow","Actual Inflow","Actual Outflow")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $wsMulti.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
    @(45688, 100000, 75000, 95000, 74000),
    @(45716, 100000, 75000, 105000, 76000),
    @(45747, 100000, 75000, 98000, 77000),
    @(45777, 100000, 75000, 100000, 78000),
    @(45808, 100000, 75000, 102000, 79000),
    @(45838, 100000, 75000, 99000, 80000),
    @(45869, 100000, 75000, 96000, 81000),
    @(45900, 100000, 75000, 97000, 82000),
    @(45930, 100000, 75000, 101000, 83000),
    @(45961, 100000, 75000, 95000, 84000),
    @(45991, 100000, 75000, 97000, 85000),
    @(46022, 100000, 75000, 98000, 86000)
)

$r = 2
foreach ($row in $rows) {
    $wsMulti.Cells.Item($r, 1).Value = $row[0]
    $wsMulti.Cells.Item($r, 1).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
    $wsMulti.Cells.Item($r, 2).Value = "Unit_A"
    $wsMulti.Cells.Item($r, 3).Value = "MS1"
    $wsMulti.Cells.Item($r, 4).Value = "USD"
    $wsMulti.Cells.Item($r, 5).Value = "USD_THB"
    $wsMulti.Cells.Item($r, 6).Value = $row[1]
    $wsMulti.Cells.Item($r, 7).Value = $row[2]
    $wsMulti.Cells.Item($r, 8).Value = $row[3]
    $wsMulti.Cells.Item($r, 9).Value = $row[4]
    $r = $r + 1
}

$wsMulti.Range("D1:D4").Select()

# ---------------------------------------------------------------------------
# 3. "As-Sold" sheet: add a "Currency" column (F) with header styling
#    (bold font, thin border all round, centered/top aligned) and "USD" in
#    every data row. The "Currency"/"USD" strings already exist in the
#    shared-string table from step 2, so no new entries are created here.
# ---------------------------------------------------------------------------
$wsAsSold = $wb.Worksheets.Item("As-Sold")
$wsAsSold.Activate()

$wsAsSold.Range("F1").Value = "Currency"
$wsAsSold.Range("F1").Font.Bold = $true
$wsAsSold.Range("F1").HorizontalAlignment = -4108
$wsAsSold.Range("F1").VerticalAlignment = -4160
$wsAsSold.Range("F1").Borders.LineStyle = 1

$wsAsSold.Range("F2").Value = "USD"
$wsAsSold.Range("F3").Value = "USD"
$wsAsSold.Range("F4").Value = "USD"
$wsAsSold.Range("F5").Value = "USD"
$wsAsSold.Range("F6").Value = "USD"
$wsAsSold.Range("F7").Value = "USD"

$wsAsSold.Range("F1:F7").Select()

# ---------------------------------------------------------------------------
# 4. "Actual" sheet: same "Currency" column, one column further right (G).
# ---------------------------------------------------------------------------
$wsActual = $wb.Worksheets.Item("Actual")
$wsActual.Activate()

$wsActual.Range("G1").Value = "Currency"
$wsActual.Range("G1").Font.Bold = $true
$wsActual.Range("G1").HorizontalAlignment = -4108
$wsActual.Range("G1").VerticalAlignment = -4160
$wsActual.Range("G1").Borders.LineStyle = 1

$wsActual.Range("G2").Value = "USD"
$wsActual.Range("G3").Value = "USD"
$wsActual.Range("G4").Value = "USD"
$wsActual.Range("G5").Value = "USD"
$wsActual.Range("G6").Value = "USD"
$wsActual.Range("G7").Value = "USD"

$wsActual.Range("G1:G7").Select()

# ---------------------------------------------------------------------------
# 5. "VO" sheet: same "Currency" column at E (only 4 data rows here).
# ---------------------------------------------------------------------------
$wsVO = $wb.Worksheets.Item("VO")
$wsVO.Activate()

$wsVO.Range("E1").Value = "Currency"
$wsVO.Range("E1").Font.Bold = $true
$wsVO.Range("E1").HorizontalAlignment = -4108
$wsVO.Range("E1").VerticalAlignment = -4160
$wsVO.Range("E1").Borders.LineStyle = 1

$wsVO.Range("E2").Value = "USD"
$wsVO.Range("E3").Value = "USD"
$wsVO.Range("E4").Value = "USD"

$wsVO.Range("E1:E4").Select()

# ---------------------------------------------------------------------------
# 6. "Variance" sheet: same "Currency" column at F, and move the lingering
#    selection to E16.
# ---------------------------------------------------------------------------
$wsVariance = $wb.Worksheets.Item("Variance")
$wsVariance.Activate()

$wsVariance.Range("F1").Value = "Currency"
$wsVariance.Range("F1").Font.Bold = $true
$wsVariance.Range("F1").HorizontalAlignment = -4108
$wsVariance.Range("F1").VerticalAlignment = -4160
$wsVariance.Range("F1").Borders.LineStyle = 1

$wsVariance.Range("F2").Value = "USD"
$wsVariance.Range("F3").Value = "USD"
$wsVariance.Range("F4").Value = "USD"
$wsVariance.Range("F5").Value = "USD"
$wsVariance.Range("F6").Value = "USD"
$wsVariance.Range("F7").Value = "USD"

$wsVariance.Range("E16").Select()

# ---------------------------------------------------------------------------
# 7. Leave "Variance" as the active tab, matching the original file.
# ---------------------------------------------------------------------------
$wsVariance.Activate()
